$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The checklist previously had two empty "to do" columns (F: Anforderungsdokument,
# G: Testreport(?)). Both the "Installationsanleitung" (col E) and the newly
# introduced "User Stories Dokument" (col F, replacing the old F/G pair) have now
# been checked, so fill in the results and retire column G.

# --- Header row: rename the last checked-document column ---
$ws.Range("F1").Value = "User Stories Dokument"
$ws.Range("G1").Value = $null

# --- Column E results: Installationsanleitung ---
$ws.Range("E2").Value = "ja, angepasst"
$ws.Range("E3").Value = "ja "
$ws.Range("E4").Value = "ja "
$ws.Range("E5").Value = "ja "
$ws.Range("E6").Value = "ja "
$ws.Range("E7").Value = "ja "
$ws.Range("E8").Value = "ja "
$ws.Range("E9").Value = "ja "
$ws.Range("E10").Value = "ja "
$ws.Range("E11").Value = "ja "
$ws.Range("E12").Value = "ja "
$ws.Range("E13").Value = "ja "
$ws.Range("E14").Value = "ja "
$ws.Range("E15").Value = "ja"
$ws.Range("E16").Value = "nein"
$ws.Range("E17").Value = "ja "
$ws.Range("E18").Value = "gecheckt am 17.06.16"

# --- Column F results: User Stories Dokument ---
$ws.Range("F2").Value = "ja, angepasst"
$ws.Range("F3").Value = "ja"
$ws.Range("F4").Value = "ja "
$ws.Range("F5").Value = "ja "
$ws.Range("F6").Value = "ja "
$ws.Range("F7").Value = "ja "
$ws.Range("F8").Value = "nicht nötig"
$ws.Range("F9").Value = "nicht nötig"
$ws.Range("F10").Value = "ja "
$ws.Range("F11").Value = "keine schachtelung nötig"
$ws.Range("F12").Value = "ja "
$ws.Range("F13").Value = "ja "
$ws.Range("F14").Value = "weder noch vorhanden"
$ws.Range("F15").Value = "nicht vorhanden"
$ws.Range("F16").Value = "nicht vorhanden"
$ws.Range("F17").Value = "nicht vorhanden"
$ws.Range("F18").Value = "gecheckt am 17.06.2016"

# Drop the now-obsolete column G (Anforderungsdokument/Testreport data never got filled in)
$ws.Columns("G").Delete()

# Columns D and E hold internal/raw check data and are no longer interesting to
# readers, so narrow/hide them and shrink column B; column C is reserved/hidden too.
$ws.Columns("B").ColumnWidth = 54.29
$ws.Columns("C").ColumnWidth = 13.14
$ws.Columns("C").Hidden = $true
$ws.Columns("D").ColumnWidth = 21.86
$ws.Columns("D").Hidden = $true
$ws.Columns("E").ColumnWidth = 20.43
$ws.Columns("E").Hidden = $true

# Reflect where the author last clicked while reviewing the sheet
$ws.Range("F19").Select()
